$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(1, 1).Value2 = 'RECALL_CLASSIFICATION_DATE'
$ws.Cells.Item(1, 2).Value2 = 'PRODUCT_TYPE'
$ws.Cells.Item(1, 3).Value2 = 'CLASSIFICATION'
$ws.Cells.Item(1, 4).Value2 = 'RECALL_NUMBER'
$ws.Cells.Item(1, 5).Value2 = 'PRODUCT'
$ws.Cells.Item(1, 6).Value2 = 'RECALLING_FIRM'
$ws.Cells.Item(1, 7).Value2 = 'MANUFACTURER'
$ws.Cells.Item(1, 8).Value2 = 'RECALL_INITIALIZATION_DATE'
$ws.Cells.Item(1, 9).Value2 = 'REASON'
$ws.Cells.Item(1, 10).Value2 = 'VOLUME'
$ws.Cells.Item(1, 11).Value2 = 'DISTRIBUTION'
$ws.Cells.Item(2, 3).Value2 = 'II'
$ws.Cells.Item(2, 5).Value2 = 'Monoject Hypodermic Needles Single Use Epoxy Needle Attachment Polypropylene Hub Luer Lock Hub 19 Ga'
$ws.Cells.Item(2, 7).Value2 = 'Sherwood Medical Company Deland Florida'
$ws.Cells.Item(2, 8).Value2 = 'December 21 1994'
$ws.Cells.Item(2, 9).Value2 = 'The Needles May Contain A Coating Of Silicone Oil In The Interior Of The Needle Hub'
$ws.Cells.Item(2, 10).Value2 = '642 600 Needles Of Lot 237593 And 276 900 Needles Of Lot 238057 Were Distributed'
$ws.Cells.Item(3, 3).Value2 = 'II'
$ws.Cells.Item(3, 5).Value2 = 'Chemspor Biological And Chemical Indicator A Bi For Monitoring Steam Sterilization Processes In Wet Environments'
$ws.Cells.Item(3, 7).Value2 = 'Amsco International Inc'
$ws.Cells.Item(3, 8).Value2 = 'February 22 1995'
$ws.Cells.Item(3, 9).Value2 = 'The Device Is Labeled For 48 Hour Incubation And 510 K Clearance Is For 7 Days'
$ws.Cells.Item(3, 10).Value2 = '9 988 Boxes Of Product Were Distributed'
$ws.Cells.Item(4, 3).Value2 = 'II'
$ws.Cells.Item(4, 5).Value2 = 'Proof Flash Self Contained Biological And Chemical Indicator A Bi For Monitoring Steam 270 Xf Steam Gravity Sterilization Cycles'
$ws.Cells.Item(4, 7).Value2 = 'Amsco International Inc'
$ws.Cells.Item(4, 8).Value2 = 'February 22 1995'
$ws.Cells.Item(4, 9).Value2 = 'The Product Is Labeled For A 24 Hour Or 48 Hour Incubation Period'
$ws.Cells.Item(4, 10).Value2 = '25 462 Cases 50 Vials Per Box Were Distributed'
$ws.Cells.Item(5, 3).Value2 = 'II'
$ws.Cells.Item(5, 5).Value2 = 'Saxon Gold Ultra Spermicidal Lubricated Latex Condoms With Nonoxynol 9 Packed In Plastic Film'
$ws.Cells.Item(5, 7).Value2 = 'Safetex Division Of Aladan Corporation Colonial Heights Virginia'
$ws.Cells.Item(5, 8).Value2 = 'September 21 1994'
$ws.Cells.Item(5, 9).Value2 = 'Condoms Were Found To Have A Lower Concentration Of The Spermicide Nonoxynol 9 Than Was Represented In The Device''s Labeling And Premarket Notification'
$ws.Cells.Item(5, 10).Value2 = '20 Gross Were Distributed Firm Estimates None Remains On The Market'
$ws.Cells.Item(5, 11).Value2 = 'New York Kansas'
$ws.Cells.Item(6, 3).Value2 = 'II'
$ws.Cells.Item(6, 5).Value2 = 'Mrinnervu Disposable Endorectal Coils For Mri Imaging For The Colon And Cervix A Cervix Coil Model Bcr 15 B Colorectal Coil Model Pcc 15'
$ws.Cells.Item(6, 7).Value2 = 'Medrad Inc'
$ws.Cells.Item(6, 8).Value2 = 'January 31 1995'
$ws.Cells.Item(6, 10).Value2 = 'Approximately A 812 Units B 915 Units Were Distributed Firm Estimated That 334 Rectal Units And 252 Cervical Units Remained On Market At Time Of Recall Initiation'
$ws.Cells.Item(7, 3).Value2 = 'II'
$ws.Cells.Item(7, 5).Value2 = 'Coveen Urisheath Uriliner Male External Catheters A Pediatric Geriatric Catheter Product 5120 H 25 Mm B Large Catheter Product 5135 H 35 Mm'
$ws.Cells.Item(7, 7).Value2 = 'Coloplast A S Humlebaek Denmark'
$ws.Cells.Item(7, 8).Value2 = 'January 31 1995'
$ws.Cells.Item(7, 9).Value2 = 'The Labeling Incorrectly Identifies The Product As Non Latex Product'
$ws.Cells.Item(7, 10).Value2 = '3 500 Units Of Lot 86811 And 200 Units Of Lot 86764 Were Distributed'
$ws.Cells.Item(8, 3).Value2 = 'II'
$ws.Cells.Item(8, 5).Value2 = 'Adult Child And Infant Cpr Resuscitator Bags Single Patient Use Disposable With Related Accessories Mask Peep Valve Peep Adapter And Manometer Port'
$ws.Cells.Item(8, 7).Value2 = 'Plaxtron Industrial M Sdn'
$ws.Cells.Item(8, 8).Value2 = 'August 29 1994'
$ws.Cells.Item(8, 9).Value2 = 'The Devices Have Excess Glue Around The Intake Valve On The Interior Of The Bag'
$ws.Cells.Item(8, 10).Value2 = '23 610 Units Were Distributed'
$ws.Cells.Item(9, 3).Value2 = 'II'
$ws.Cells.Item(9, 5).Value2 = 'Transflator Ii Balloon Inflation Syringes A Transflator Balloon Inflation Syringe Catalog 05 0802 B Transflator Ii Balloon Syringe With Ptca Accessory Kit Catalog 05 0804'
$ws.Cells.Item(9, 7).Value2 = 'Surgical Technologies Inc'
$ws.Cells.Item(9, 8).Value2 = 'February 2 1995'
$ws.Cells.Item(9, 9).Value2 = 'Devices Have Broken Transmitting Led Solder Joints Which May Result In Loss Of Transmission And The Signal Lost Read Out On The Monitors'
$ws.Cells.Item(9, 10).Value2 = '1 732 Units Were Distributed Firm Estimates None Remains On The Market'
$ws.Cells.Item(9, 11).Value2 = 'Nationwide Belgium Canada Italy The Netherlands United Kingdom'
$ws.Cells.Item(10, 3).Value2 = 'II'
$ws.Cells.Item(10, 5).Value2 = 'Medical Research Laboratory Dpd Defibrillator Sold Under Datascope Part 0998 00 0004 Used In Conjunction With Datascope Passport Patient Monitoring System'
$ws.Cells.Item(10, 7).Value2 = 'Medical Research Laboratories Buffalo Grove Illinois'
$ws.Cells.Item(10, 8).Value2 = 'March 2 1995'
$ws.Cells.Item(10, 9).Value2 = 'There Is A Potential Of The Devices When Used In The Synchronous Mode With The Datascope Passport Monitor For The Purpose Of Cardioversion To Discharge As Late As Approximately 160 Ms After The Peak Of The R Wave Of The Electrocardiogram'
$ws.Cells.Item(10, 10).Value2 = '656 Units Were Distributed'
$ws.Cells.Item(11, 3).Value2 = 'II'
$ws.Cells.Item(11, 5).Value2 = 'Encore Orthopedic Foundation Brand Knee Tibial Insert 13 Mm Size 4'
$ws.Cells.Item(11, 7).Value2 = 'Encore Orthopedics Inc'
$ws.Cells.Item(11, 8).Value2 = 'October 7 1994'
$ws.Cells.Item(11, 9).Value2 = 'The Label On The Interior Pouch May Incorrectly Identify The Insert Thickness As 15 Mm However The Actual Thickness Is 13 Mm Which Is Correctly Identified On The Exterior Box Label'
$ws.Cells.Item(11, 10).Value2 = '15 Inserts Were Distributed Firm Estimates None Remain On Market'
$ws.Cells.Item(11, 11).Value2 = 'New York Oklahoma Pennsylvania South Carolina Utah Texas Washington State'
$ws.Cells.Item(12, 3).Value2 = 'II'
$ws.Cells.Item(12, 5).Value2 = 'Slit Rotation Knob On The Retroskop Cs An Illumination Module Microscope Attachment'
$ws.Cells.Item(12, 7).Value2 = 'Carl Zeiss Inc'
$ws.Cells.Item(12, 8).Value2 = 'February 6 1995'
$ws.Cells.Item(12, 9).Value2 = 'The Knob Originally Designed To Be Removable So It Could Be Sterilized Has The Potential To Loosen And Fall Off From The Retroskop After Being Sterilized'
$ws.Cells.Item(12, 10).Value2 = '80 Attachments Were Distributed'
$ws.Cells.Item(12, 11).Value2 = 'Nationwide Canada'
$ws.Cells.Item(13, 3).Value2 = 'II'
$ws.Cells.Item(13, 5).Value2 = 'D 4600 Air Exchange Line Sterile Single Patient Use An Accessory For Ophthalmic Use'
$ws.Cells.Item(13, 7).Value2 = 'Storz Instrument Company St'
$ws.Cells.Item(13, 8).Value2 = 'November 28 1994'
$ws.Cells.Item(13, 9).Value2 = 'Due To A Packaging Deficiency The Pouches Were Improperly Sealed Therefore Compromising Sterility'
$ws.Cells.Item(13, 10).Value2 = '2 213 Boxes 10 Pouches Per Box Were Distributed'
$dcell = $ws.Cells.Item(14, 1)
$dcell.Formula = '="1995-03-29"'
$dcell.Copy()
$dcell.PasteSpecial(-4163)
$ws.Cells.Item(14, 2).Value2 = 'Devices'
$ws.Cells.Item(14, 3).Value2 = 'III'
$ws.Cells.Item(14, 4).Value2 = 'z-442/443-5'
$ws.Cells.Item(14, 5).Value2 = 'Mmi Barium Enema Kits For Use In X Ray Diagnosis Of The Lower Gastrointestinal Tract A Catalog 354 00 Assembled Empty Barium Enema Kit B Catalog 381 00 Assembled Filled Barium Enema Kit'
$ws.Cells.Item(14, 6).Value2 = 'Mmi'
$ws.Cells.Item(14, 7).Value2 = 'Ridge Products Lakewood New Jersey'
$ws.Cells.Item(14, 8).Value2 = 'May 20 1994'
$ws.Cells.Item(14, 9).Value2 = 'Devices May Have A Fracture In The Seal At The Base Of The Bag Which Could Result In Leakage Product Mdt Mckesson Desert Air Plus Dental Air Compressors A Model 301862 Single Head 115 V B Model 301863 Single Head 230 V C Model 301864 Twin Head 230 V D Model 301865 Triple Head 230 V Recall Z 502 505 5 Code Serial Numbers A An 0001 Kn 0122 And A 00127 B Bn 0001 Gn 0064 C An 0001 Kn 0164 And B 00165 D Cn 0001 Kn 0015 And Ln 0016 And Ln 0017 Manufacturer Mdt Diagnostic Company North Charleston South Carolina Recalled By Manufacturer By Visit Beginning June 10 1994 Firm Initiated Field Correction Complete Distribution Nationwide And Canada Quantity A 124 Units B 65 Units C 166 Units D 18 Units Were Distributed Reason The Nylon Coalescing Filter Bowl May Break Apart Due To Vibrational Fatigue Product Mdt Mckesson Desert Air Plus Dental Air Compressors A Model 301862 Single Head 115 V B Model 301863 Single Head 230 V C Model 301864 Twin Head 230 V D Model 301865 Triple Head 230 V Recall Z 506 509 5 Code Serial Numbers A An 0001 Kn 0122 B Bn 0001 Gn 0064 C An 0001 Kn 0164 D Cn 0001 Kn 0015 Manufacturer Mdt Diagnostic Company North Charleston South Carolina Recalled By Manufacturer By Letter December 6 1993 Firm Initiated Field Correction Complete Distribution Nationwide And Canada Quantity A 123 Units B 65 Units C 165 Units D 16 Units Were Distributed Reason The Clear Styrene Desiccant Enclosure May Explode Product Bausch Lomb Gold Medalist Toric Visibility Tinted Contact Lenses Labeled With The Following Parameters A Blister Lot Number Sy 009541 Bc 8 6 Sph 2 00 Cyl 125 Ax 90 B Blister Lot Number Sy 009430 Bc 8 3 Sph 4 00 Cyl 1 75 Ax 170 Recall Z 551 552 5 Code Carton Lot Numbers As Follows A Sz 064942 Through Sz 064956 Sz 069065 Sz 065693 8 Sz 066011 Sz 065748 Sz 066017 Sz 066019 Sz 065751 Sz 066022 Sz 066263 Sz 065749 Sz 065752 Sz 065757 Sz 065759 B Sz 064968 Sz 064972 Sz 064973 Sz 064974 Sz 064977 Sz 064978 Sz 065265 Sz 065490 Sz 065231 Sz 065232 Sz 065648 Sz 065549 Sz 065267 Sz 065271 Sz 065490 Manufacturer Bausch Lomb Inc Sarasota Florida Recalled By Bausch Lomb Inc Rochester New York By Telephone And Fax On January 24 And 25 1995 Firm Initiated Recall Complete Distribution Alaska Kansas Minnesota Virginia California Florida Georgia Illinois Iowa Missouri New Jersey New York Ohio Oregon Wisconsin Canada Australia Ireland Quantity A 66 Lenses B 71 Lenses Were Distributed Reason An Unidentified Number Of Lenses Are Labeled With The Incorrect Base Curve Sphere Cylinder And Axis Parameters Product Dade Alcohol And Ammonia Control Level I Normal And Level Ii High Recall Z 560 5 Code Catalog No Lot No Exp Date B 5118 Am 4 Pac 34 2 8 95 B 5188 Am 5 Pac 35 10 9 95 B 5188 S Pac 34 35 2 8 10 9 95 B 5118 Pac 36 12 25 95 Manufacturer Baxter Diagnostics Inc Miami Florida Recalled By Manufacturer By Letter August 1994 Firm Initiated Field Correction Complete Distribution Nationwide And International Quantity 29 022 Units Were Distributed Reason The Labeling Was Misleading As To The Ranges Listed In The Package Insert Product Dade Immunoassay Controls Tri Level Level I Level Ii And Level Iii A Catalog B 5700 5 B Catalog B 5700 6 C Catalog B 5700 7 D Catalog B 5700 8 Recall Z 561 564 5 Code Lot Numbers A Iack 22 M B Iac 1 112 M C Iac 2 212 M D Iac 3 312 M Manufacturer Baxter Diagnostics Inc Miami Florida Recalled By Manufacturer By Letter On Or About October 20 1994 Firm Initiated Field Correction Complete Distribution Texas Michigan California Wisconsin Florida New York Quantity 35 Kits Were Distributed Reason The Labeling Is Misleading As The Results For The Free T 4 Clinical Assays 2 Step Cal 535 A And Direct Cal 535 S Were Reversed On The Clinical Assay Sheets Provided With The Test Kits 10 Seizures Product Frozen Scallops 95 710 318 Charge Adulterated The Article Consists In Whole Or In Part Of Decomposed Scallops Location Continental Freezer Of Illinois Chicago Illinois Filed March 22 1995 U S District Court For The Northern District Of Illinois Civ 95 C 1767 Fdc 67074 Seized March 22 1995 Goods Valued At Approximately 9 320 Product Frozen Peeled Shrimp And Imitation Breaded Shrimp 94 681 944 5 Charge Adulterated The Articles Consist In Whole Or In Part Of A Decomposed Substance By Reason Of The Presence Therein Of Decomposed Shrimp Location Seaboard Cold Storage Inc Tampa Florida Filed January 30 1995 U S District Court For The Middle District Of Florida Tampa Division Civ 95 129 Civ T 23 B Fdc 67053 Seized March 9 1995 Goods Valued At Approximately 9 500 10 End Of Enforcement Report For March 29 1995 Blank Pages May Follow'
$ws.Cells.Item(14, 10).Value2 = 'A 242 Kits B 48 Kits Were Distributed'
$ws.Cells.Item(14, 11).Value2 = 'New York Florida Michigan'

$excel.CutCopyMode = 0
Write-Output "applied"
